$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Grupo: ... Sistema: ... Fecha: ..." paragraph:
#    split the single run into several runs, fill in the group number ("4")
#    and the system name ("Optical Marketing", with "Optical" flagged by the
#    spell-checker via <w:proofErr>), and shorten the blank filler before
#    "Sistema:"/"Fecha:" accordingly.
# ---------------------------------------------------------------------------
$grupoPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Grupo:*Sistema:*Fecha:*") {
        $grupoPara = $p
        break
    }
}

if ($grupoPara -ne $null) {
    # Range covering just the paragraph's text (excluding the paragraph mark).
    $pRange = $grupoPara.Range
    $textRange = $d.Range($pRange.Start, $pRange.End - 1)

    $xmlFrag = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">Grupo:   </w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>4</w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">                                             Sistema: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>Optical</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> Marketing</w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">                                           Fecha:                   </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $textRange.InsertXML($xmlFrag)
}

# ---------------------------------------------------------------------------
# 2) "Sprint: ... Semana: ..." paragraph: drop one space between "Sprint:"
#    and "Semana:".
# ---------------------------------------------------------------------------
$sprintPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Sprint:*Semana:*") {
        $sprintPara = $p
        break
    }
}

if ($sprintPara -ne $null) {
    $pRange2 = $sprintPara.Range
    $textRange2 = $d.Range($pRange2.Start, $pRange2.End - 1)
    $textRange2.Text = "Sprint:                                                  Semana:               "
}
